$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.151.29'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '3.399.00'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.09'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.398.69'
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('E10').Value = '  +2.86%  '
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.381'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('D13').Value = '3.977.74'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.66'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000172'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').Value = '3.397.81'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').Value = '61.174.76'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.52%  '
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '3.533.58'
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.553'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('E27').Value = '  -2.65%  '
$ws.Range('E28').Value = '  +9.38%  '
$ws.Range('E29').Value = '  -6.61%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.46'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E36').Value = '  +2.04%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.90'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.12'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '166.24'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '25.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.75'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.16%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('E47').Value = '  -2.62%  '
$ws.Range('D48').Value = '2.524.53'
$ws.Range('E48').Value = '  +6.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.48%  '
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('E51').Value = '  -0.23%  '
